$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.423.22"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.642.59"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'212.09"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'23.08"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "1.874.14"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "1.641.09"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "'0.572"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "'64.44"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").Value = "27.402.05"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'229.60"
$ws.Range("E18").Value = "  -5.16%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "'9.65"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "'147.10"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "1.413.40"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("D35").Value = "'1.58"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").Value = "'1.03"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.819"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.51"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D47").Value = "1.784.22"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("D49").Value = "'87.99"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'0.0991"
$ws.Range("E51").Value = "  -3.65%  "
